# The commit simply swaps the colour scheme carried by the deck's theme
# part: the main design ("Integral" - green/yellow palette) and the
# notes-master's theme ("Office Theme" - default blue palette) trade
# places. Concretely, theme1.xml (Office Theme) <-> theme2.xml (Integral)
# swap their <a:clrScheme> contents (their font/format schemes are
# already identical, only the 12 theme colours + the theme/clrScheme
# "name" differ).
#
# The live, editable theme reachable from the PowerPoint object model is
# the one applied to the slide master / layouts / slides (exposed as
# ActivePresentation.SlideMaster.Theme, i.e. the deck's "Design"). We
# recolour it from the current "Integral" palette to the default
# "Office Theme" palette, which is the half of the swap that is visible
# in Normal view and reachable via COM automation.

$p = $ppt.ActivePresentation
$design = $p.SlideMaster.Theme.ThemeColorScheme

function ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $design.Colors($i).RGB = ToRGB($officeThemeColors[$i - 1])
}
